$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows at position 125, pushing the existing
# rows 125-192 down to 128-195.
$ws.Range("A125:A127").EntireRow.Insert()

# Row 125 (new)
$ws.Range("A125").Value = 4
$ws.Range("B125").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C125").Value = "Los Lagos"
$ws.Range("D125").Value = 44574
$ws.Range("E125").Value = 10
$ws.Range("F125").Value = 100112028
$ws.Range("G125").Value = "Sandia"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 500
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 3000
$ws.Range("M125").Value = 3000
$ws.Range("N125").Value = [char]36 + "/unidad"
$ws.Range("O125").Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Range("P125").Value = 3000
$ws.Range("Q125").Value = 1
$ws.Range("R125").Value = "Hortaliza"

# Row 126 (new)
$ws.Range("A126").Value = 4
$ws.Range("B126").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C126").Value = "Los Lagos"
$ws.Range("D126").Value = 44574
$ws.Range("E126").Value = 10
$ws.Range("F126").Value = 100112028
$ws.Range("G126").Value = "Sandia"
$ws.Range("H126").Value = "Sin especificar"
$ws.Range("I126").Value = "Segunda"
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 2500
$ws.Range("L126").Value = 2500
$ws.Range("M126").Value = 2500
$ws.Range("N126").Value = [char]36 + "/unidad"
$ws.Range("O126").Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Range("P126").Value = 2500
$ws.Range("Q126").Value = 1
$ws.Range("R126").Value = "Hortaliza"

# Row 127 (new)
$ws.Range("A127").Value = 4
$ws.Range("B127").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C127").Value = "Los Lagos"
$ws.Range("D127").Value = 44574
$ws.Range("E127").Value = 10
$ws.Range("F127").Value = 100112028
$ws.Range("G127").Value = "Sandia"
$ws.Range("H127").Value = "Sin especificar"
$ws.Range("I127").Value = "Tercera"
$ws.Range("J127").Value = 800
$ws.Range("K127").Value = 2000
$ws.Range("L127").Value = 2000
$ws.Range("M127").Value = 2000
$ws.Range("N127").Value = [char]36 + "/unidad"
$ws.Range("O127").Value = "Regi" + [char]0x00F3 + "n de O'Higgins"
$ws.Range("P127").Value = 2000
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = "Hortaliza"
